{"js": "// Update the date line and the 25 multiplication problems in the table\n// to match the new day's worksheet content.\nconst replacements = [\n  { old: \"2023-12-03 Sunday\", new: \"2023-12-04 Monday\" },\n  { old: \"53\u00d745=\", new: \"72\u00d729=\" },\n  { old: \"57\u00d793=\", new: \"93\u00d758=\" },\n  { old: \"48\u00d742=\", new: \"79\u00d713=\" },\n  { old: \"36\u00d783=\", new: \"76\u00d772=\" },\n  { old: \"37\u00d730=\", new: \"80\u00d750=\" },\n  { old: \"91\u00d775=\", new: \"43\u00d769=\" },\n  { old: \"72\u00d730=\", new: \"58\u00d772=\" },\n  { old: \"63\u00d796=\", new: \"47\u00d757=\" },\n  { old: \"98\u00d781=\", new: \"36\u00d783=\" },\n  { old: \"18\u00d756=\", new: \"71\u00d743=\" },\n  { old: \"80\u00d791=\", new: \"45\u00d711=\" },\n  { old: \"48\u00d776=\", new: \"69\u00d772=\" },\n  { old: \"71\u00d758=\", new: \"96\u00d719=\" },\n  { old: \"69\u00d799=\", new: \"91\u00d717=\" },\n  { old: \"96\u00d764=\", new: \"88\u00d726=\" },\n  { old: \"58\u00d791=\", new: \"34\u00d758=\" },\n  { old: \"15\u00d716=\", new: \"97\u00d713=\" },\n  { old: \"17\u00d776=\", new: \"99\u00d789=\" },\n  { old: \"67\u00d744=\", new: \"98\u00d793=\" },\n  { old: \"46\u00d771=\", new: \"99\u00d777=\" },\n  { old: \"44\u00d774=\", new: \"48\u00d759=\" },\n  { old: \"21\u00d787=\", new: \"19\u00d784=\" },\n  { old: \"29\u00d750=\", new: \"42\u00d797=\" },\n  { old: \"17\u00d736=\", new: \"55\u00d791=\" },\n  { old: \"63\u00d775=\", new: \"59\u00d729=\" },\n];\n\nconst body = context.document.body;\n\n// Apply each replacement one at a time: search the ORIGINAL text (unique in\n// the document) and replace just that hit. Doing this sequentially (rather\n// than one big search pass) avoids any risk of a newly-inserted value being\n// re-matched by a later search (e.g. \"98\u00d781=\" becomes \"36\u00d783=\", which was\n// itself a pre-existing original value elsewhere in the table).\nfor (const { old, new: replacement } of replacements) {\n  const results = body.search(old, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${old}`);\n  }\n  results.items[0].insertText(replacement, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 multiplication problems in the table\n# to match the new day's worksheet content.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = '2023-12-03 Sunday'; New = '2023-12-04 Monday' },\n    @{ Old = '53\u00d745='; New = '72\u00d729=' },\n    @{ Old = '57\u00d793='; New = '93\u00d758=' },\n    @{ Old = '48\u00d742='; New = '79\u00d713=' },\n    @{ Old = '36\u00d783='; New = '76\u00d772=' },\n    @{ Old = '37\u00d730='; New = '80\u00d750=' },\n    @{ Old = '91\u00d775='; New = '43\u00d769=' },\n    @{ Old = '72\u00d730='; New = '58\u00d772=' },\n    @{ Old = '63\u00d796='; New = '47\u00d757=' },\n    @{ Old = '98\u00d781='; New = '36\u00d783=' },\n    @{ Old = '18\u00d756='; New = '71\u00d743=' },\n    @{ Old = '80\u00d791='; New = '45\u00d711=' },\n    @{ Old = '48\u00d776='; New = '69\u00d772=' },\n    @{ Old = '71\u00d758='; New = '96\u00d719=' },\n    @{ Old = '69\u00d799='; New = '91\u00d717=' },\n    @{ Old = '96\u00d764='; New = '88\u00d726=' },\n    @{ Old = '58\u00d791='; New = '34\u00d758=' },\n    @{ Old = '15\u00d716='; New = '97\u00d713=' },\n    @{ Old = '17\u00d776='; New = '99\u00d789=' },\n    @{ Old = '67\u00d744='; New = '98\u00d793=' },\n    @{ Old = '46\u00d771='; New = '99\u00d777=' },\n    @{ Old = '44\u00d774='; New = '48\u00d759=' },\n    @{ Old = '21\u00d787='; New = '19\u00d784=' },\n    @{ Old = '29\u00d750='; New = '42\u00d797=' },\n    @{ Old = '17\u00d736='; New = '55\u00d791=' },\n    @{ Old = '63\u00d775='; New = '59\u00d729=' }\n)\n\n# Apply each replacement one at a time against a fresh Range over the whole\n# document: search for the ORIGINAL text (unique in the document) and swap in\n# the new value. wdReplaceAll (2) here only ever touches the single unique\n# match for that search string. Processing sequentially in document order\n# means no later search text collides with a value just written earlier\n# (e.g. \"98\u00d781=\" becomes \"36\u00d783=\", which was itself a pre-existing original\n# value elsewhere in the table but is only searched for BEFORE that point).\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n    $found = $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $found) {\n        throw \"Text not found: $($r.Old)\"\n    }\n}\n"}
